$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header values in row 1: P1=14, Q1=15, matching the style of the
# existing header cells (B1:O1 use style index 1 -> copy from O1).
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats

# Swap the I/K and M/O column values for data rows 2-25.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1

    # New columns P and Q, value 2 for every data row.
    $ws.Cells.Item($r, 16).Value = 2
    $ws.Cells.Item($r, 17).Value = 2
}
